$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Teun Jaspers" (row 12) and "Henri Vogels" (row 13) rows as
# Whitelist = TRUE, copying the existing D15 cell (already a "TRUE" text
# entry) so the new cells pick up the exact same shared-string value/type
# and formatting used elsewhere in column D (D15, D26) instead of getting
# auto-converted to a boolean literal.
$ws.Range("D15").Copy($ws.Range("D12"))
$ws.Range("D15").Copy($ws.Range("D13"))

# Move the active selection/scroll position to D13, where the second new
# Whitelist flag was just entered.
[void]$ws.Range("D13").Select()
